# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 15 de Julio de 2020 a las 02:36"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 3544857
$ws.Range("C4").Value = 65374
$ws.Range("D4").Value = 1597066
$ws.Range("E4").Value = 1808654
$ws.Range("G4").Value = 929
$ws.Range("H4").Value = 139137

# Row 19 - Alemania
$ws.Range("B19").Value = 200766
$ws.Range("C19").Value = 330
$ws.Range("E19").Value = 6122
$ws.Range("G19").Value = 5
$ws.Range("H19").Value = 9144

# Row 145 - Uruguay
$ws.Range("B145").Value = 997
$ws.Range("C145").Value = 8
$ws.Range("D145").Value = 905
$ws.Range("E145").Value = 61

# Row 166 - Comoras
$ws.Range("B166").Value = 321
$ws.Range("C166").Value = 4
$ws.Range("D166").Value = 302
$ws.Range("E166").Value = 12

# Row 167 - Guyana
$ws.Range("B167").Value = 308
$ws.Range("C167").Value = 8
$ws.Range("D167").Value = 156
$ws.Range("E167").Value = 135

# Row 190 - Islas Turcas y Caicos
$ws.Range("D190").Value = 12
$ws.Range("E190").Value = 58

# Row 196 - Belice
$ws.Range("B196").Value = 39
$ws.Range("C196").Value = 2
$ws.Range("D196").Value = 21
$ws.Range("E196").Value = 16
